# Updates cryptos list values per "Updated cryptos list" commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.057.55"
$ws.Range("E2").Value = "  +0.43%  "

$ws.Range("D3").Value = "1.923.88"
$ws.Range("E3").Value = "  +1.00%  "

$ws.Range("D4").Value = "'1.004"
$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").Value = "'325.76"
$ws.Range("E5").Value = "  +0.35%  "

$ws.Range("D6").Value = "'1.003"
$ws.Range("E6").Value = "  +0.04%  "

$ws.Range("D7").Value = "'0.4584"
$ws.Range("E7").Value = "  -0.28%  "

$ws.Range("D8").Value = "'0.3819"
$ws.Range("E8").Value = "  +0.05%  "

$ws.Range("D9").Value = "'0.07756"
$ws.Range("E9").Value = "  +0.33%  "

$ws.Range("D10").Value = "'0.9798"
$ws.Range("E10").Value = "  -0.29%  "

$ws.Range("D11").Value = "'22.63"
$ws.Range("E11").Value = "  +2.52%  "

$ws.Range("D12").Value = "1.947.10"
$ws.Range("E12").Value = "  +1.59%  "

$ws.Range("D13").Value = "'5.708"
$ws.Range("E13").Value = "  +0.42%  "

$ws.Range("D14").Value = "'6.964"
$ws.Range("E14").Value = "  -0.49%  "

$ws.Range("D15").Value = "'0.07009"
$ws.Range("E15").Value = "  -0.36%  "

$ws.Range("B16").Value = "Litecoin"
$ws.Range("C16").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D16").Value = "'84.72"
$ws.Range("E16").Value = "  +0.59%  "

$ws.Range("B17").Value = "BinanceUSD"
$ws.Range("C17").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D17").Value = "'1.006"
$ws.Range("E17").Value = "  +0.08%  "

$ws.Range("D18").Value = "'0.000009496"
$ws.Range("E18").Value = "  -0.55%  "

$ws.Range("E19").Value = "  -0.35%  "

$ws.Range("D20").Value = "'1.003"
$ws.Range("E20").Value = "  -0.01%  "

$ws.Range("D21").Value = "29.064.76"
$ws.Range("E21").Value = "  +0.40%  "

$ws.Range("D22").Value = "'5.348"
$ws.Range("E22").Value = "  +0.30%  "

$ws.Range("E23").Value = "  +0.70%  "

$ws.Range("D24").Value = "2.125.94"
$ws.Range("E24").Value = "  -1.69%  "

$ws.Range("D25").Value = "'2.057"
$ws.Range("E25").Value = "  -0.89%  "

$ws.Range("D26").Value = "'157.89"
$ws.Range("E26").Value = "  +0.80%  "

$ws.Range("D27").Value = "'18.99"
$ws.Range("E27").Value = "  -0.98%  "

$ws.Range("D28").Value = "'5.603"
$ws.Range("E28").Value = "  -0.04%  "

$ws.Range("D29").Value = "'117.59"
$ws.Range("E29").Value = "  -0.24%  "

$ws.Range("D30").Value = "'1.834"
$ws.Range("E30").Value = "  +0.25%  "

$ws.Range("D31").Value = "'0.09302"
$ws.Range("E31").Value = "  +0.44%  "

$ws.Range("D32").Value = "'0.8593"
$ws.Range("E32").Value = "  -0.18%  "

$ws.Range("D33").Value = "'5.093"
$ws.Range("E33").Value = "  -0.40%  "

$ws.Range("D34").Value = "'1.243"
$ws.Range("E34").Value = "  -0.95%  "

$ws.Range("D35").Value = "'3.012"
$ws.Range("E35").Value = "  -0.23%  "

$ws.Range("E36").Value = "  -0.65%  "

$ws.Range("D37").Value = "'1.151"
$ws.Range("E37").Value = "  +0.56%  "

$ws.Range("E38").Value = "  +0.07%  "

$ws.Range("D39").Value = "'0.02041"
$ws.Range("E39").Value = "  +0.11%  "

$ws.Range("D40").Value = "'3.093"
$ws.Range("E40").Value = "  +12.10%  "

$ws.Range("B41").Value = "TheSandbox"
$ws.Range("C41").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D41").Value = "'0.5498"
$ws.Range("E41").Value = "  -0.66%  "

$ws.Range("B42").Value = "FraxShare"
$ws.Range("C42").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D42").Value = "'7.416"
$ws.Range("E42").Value = "  -0.97%  "

$ws.Range("D43").Value = "'0.1753"
$ws.Range("E43").Value = "  -0.18%  "

$ws.Range("D44").Value = "'9.368"
$ws.Range("E44").Value = "  +0.67%  "

$ws.Range("D45").Value = "'0.000002847"
$ws.Range("E45").Value = "  +7.60%  "

$ws.Range("D46").Value = "'2.181"
$ws.Range("E46").Value = "  +3.87%  "

$ws.Range("D47").Value = "'0.5172"
$ws.Range("E47").Value = "  -0.92%  "

$ws.Range("D48").Value = "'0.06933"
$ws.Range("E48").Value = "  +1.75%  "

$ws.Range("E49").Value = "  -1.48%  "

$ws.Range("D50").Value = "'110.44"
$ws.Range("E50").Value = "  -1.39%  "

$ws.Range("D51").Value = "'1.759"
$ws.Range("E51").Value = "  -1.13%  "
